$wb = $excel.ActiveWorkbook

# This script applies a batch update of market-price-derived columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, as produced by
# the scheduled market-data refresh runner.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 35.444443
$ws.Range("I2").Value = 35.444443
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 35.444443
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 77.55555699999999
$ws.Range("N2").ClearContents()
$ws.Range("H19").Value = 1166.6666
$ws.Range("J19").Value = 1187.375
$ws.Range("L19").Value = 1187.375
$ws.Range("N19").Value = -1537.375
$ws.Range("H75").Value = 36157
$ws.Range("J75").Value = 36157
$ws.Range("L75").Value = 36157
$ws.Range("N75").Value = -38029
$ws.Range("H78").Value = 36157
$ws.Range("J78").Value = 36157
$ws.Range("L78").Value = 108471
$ws.Range("N78").Value = -117831
$ws.Range("H132").Value = 3171.111
$ws.Range("J132").Value = 2888.25
$ws.Range("L132").Value = 8664.75
$ws.Range("N132").Value = -13724.75
$ws.Range("H138").Value = 3703.383
$ws.Range("I138").Value = 869.4545000000001
$ws.Range("J138").Value = 6197.24
$ws.Range("K138").Value = 2608.3635
$ws.Range("L138").Value = 18591.72
$ws.Range("M138").Value = 2531.6365
$ws.Range("N138").Value = -28871.72

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1867664.5
$ws.Range("I32").Value = 1924970.2
$ws.Range("K32").Value = 1924970.2
$ws.Range("M32").Value = -1924683.2
$ws.Range("H61").Value = 2322.1765
$ws.Range("I61").Value = 2351.7188
$ws.Range("J61").Value = 1849.5
$ws.Range("K61").Value = 2351.7188
$ws.Range("L61").Value = 1849.5
$ws.Range("M61").Value = -2139.7188
$ws.Range("N61").Value = -2273.5
$ws.Range("H74").Value = 50218.03
$ws.Range("I74").Value = 72284.95
$ws.Range("J74").Value = 4078.0908
$ws.Range("K74").Value = 72284.95
$ws.Range("L74").Value = 4078.0908
$ws.Range("M74").Value = -71410.95
$ws.Range("N74").Value = -5826.0908
$ws.Range("H77").Value = 50218.03
$ws.Range("I77").Value = 72284.95
$ws.Range("J77").Value = 4078.0908
$ws.Range("K77").Value = 361424.75
$ws.Range("L77").Value = 20390.454
$ws.Range("M77").Value = -357056.75
$ws.Range("N77").Value = -29126.454
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H82").Value = 19999.4
$ws.Range("J82").Value = 19999.4
$ws.Range("L82").Value = 19999.4
$ws.Range("N82").Value = -20721.4
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H85").Value = 19999.4
$ws.Range("J85").Value = 19999.4
$ws.Range("L85").Value = 19999.4
$ws.Range("N85").Value = -22495.4
$ws.Range("H132").Value = 10303.5
$ws.Range("I132").Value = 10184.909
$ws.Range("K132").Value = 30554.727
$ws.Range("M132").Value = -28024.727
$ws.Range("H136").Value = 2322.1765
$ws.Range("I136").Value = 2351.7188
$ws.Range("J136").Value = 1849.5
$ws.Range("K136").Value = 7055.1564
$ws.Range("L136").Value = 5548.5
$ws.Range("M136").Value = -4505.1564
$ws.Range("N136").Value = -10648.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3375.15
$ws.Range("I105").Value = 2952.8462
$ws.Range("K105").Value = 2952.8462
$ws.Range("M105").Value = -1205.8462
$ws.Range("H134").Value = 8728.444
$ws.Range("I134").Value = 4715.4
$ws.Range("K134").Value = 14146.2
$ws.Range("M134").Value = -11611.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("H99").Value = 8469.941000000001
$ws.Range("I99").Value = 8776.556
$ws.Range("K99").Value = 8776.556
$ws.Range("M99").Value = -7278.556
$ws.Range("H105").Value = 3759.0625
$ws.Range("I105").Value = 1558.6364
$ws.Range("K105").Value = 1558.6364
$ws.Range("M105").Value = 188.3635999999999
$ws.Range("H126").Value = 8469.941000000001
$ws.Range("I126").Value = 8776.556
$ws.Range("K126").Value = 26329.668
$ws.Range("M126").Value = -23859.668
$ws.Range("H132").Value = 7153.609
$ws.Range("J132").Value = 9573.866
$ws.Range("L132").Value = 28721.598
$ws.Range("N132").Value = -33781.598
$ws.Range("H134").Value = 8951.378000000001
$ws.Range("I134").Value = 9781.157999999999
$ws.Range("K134").Value = 29343.474
$ws.Range("M134").Value = -26808.474

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 267557.94
$ws.Range("J2").Value = 668384
$ws.Range("L2").Value = 4010304
$ws.Range("N2").Value = -4010530
$ws.Range("H39").Value = 6945.385
$ws.Range("J39").Value = 8149
$ws.Range("L39").Value = 24447
$ws.Range("N39").Value = -25035
$ws.Range("H55").Value = 24080508
$ws.Range("I55").Value = 55556524
$ws.Range("J55").Value = 8342499.5
$ws.Range("K55").Value = 166669572
$ws.Range("L55").Value = 25027498.5
$ws.Range("M55").Value = -166669395
$ws.Range("N55").Value = -25027852.5
$ws.Range("H131").Value = 1809.8096
$ws.Range("J131").Value = 2050.1
$ws.Range("L131").Value = 6150.299999999999
$ws.Range("N131").Value = -16230.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 50000
$ws.Range("J86").Value = 50000
$ws.Range("L86").Value = 50000
$ws.Range("N86").Value = -52372
$ws.Range("H89").Value = 50000
$ws.Range("J89").Value = 50000
$ws.Range("L89").Value = 150000
$ws.Range("N89").Value = -161856
$ws.Range("H138").Value = 82429
$ws.Range("J138").Value = 82429
$ws.Range("L138").Value = 82429
$ws.Range("N138").Value = -92709

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3694.3333
$ws.Range("I122").Value = 2944.742
$ws.Range("K122").Value = 8834.226000000001
$ws.Range("M122").Value = -6384.226000000001
$ws.Range("H133").Value = 98000
$ws.Range("J133").Value = 98000
$ws.Range("L133").Value = 98000
$ws.Range("N133").Value = -103060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 120452.55
$ws.Range("I15").Value = 120452.55
$ws.Range("K15").Value = 120452.55
$ws.Range("M15").Value = -120164.55
$ws.Range("H136").Value = 23282772
$ws.Range("J136").Value = 44143.5
$ws.Range("L136").Value = 132430.5
$ws.Range("N136").Value = -137530.5

